# Update the "Price" (column D) and "Volume(1h)" (column E) figures for the
# cryptos list, as refreshed by the scheduled GitHub Actions job.
# Note: a leading '' (single literal apostrophe) is used for D-column values
# that look like plain numbers, so Excel stores/keeps them as text (matching
# the original inlineStr cells) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.142.02'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '3.327.98'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''582.75'
$ws.Range("E5").Value = '  +3.07%  '
$ws.Range("D6").Value = '''185.35'
$ws.Range("E6").Value = '  -3.14%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '3.323.22'
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("D9").Value = '''0.576'
$ws.Range("E9").Value = '  -2.57%  '
$ws.Range("E10").Value = '  -3.14%  '
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("D12").Value = '''47.12'
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("D14").Value = '''674.20'
$ws.Range("E14").Value = '  +9.93%  '
$ws.Range("D15").Value = '3.858.16'
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").Value = '''8.49'
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").Value = '66.244.09'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("D20").Value = '3.326.61'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '''11.10'
$ws.Range("E21").Value = '  -1.08%  '
$ws.Range("D22").Value = '''0.897'
$ws.Range("E22").Value = '  -2.41%  '
$ws.Range("D23").Value = '''17.78'
$ws.Range("E23").Value = '  -4.05%  '
$ws.Range("D24").Value = '''103.34'
$ws.Range("E24").Value = '  +1.49%  '
$ws.Range("D25").Value = '''5.05'
$ws.Range("E25").Value = '  -2.20%  '
$ws.Range("D26").Value = '''3.97'
$ws.Range("E26").Value = '  -1.41%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  -3.04%  '
$ws.Range("D29").Value = '''32.07'
$ws.Range("E29").Value = '  +4.71%  '
$ws.Range("E30").Value = '  -2.92%  '
$ws.Range("D31").Value = '''6.79'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").Value = '''599.90'
$ws.Range("E32").Value = '  +5.64%  '
$ws.Range("E33").Value = '  -5.31%  '
$ws.Range("E34").Value = '  -1.79%  '
$ws.Range("E35").Value = '  -1.24%  '
$ws.Range("D36").Value = '3.821.69'
$ws.Range("E36").Value = '  +1.86%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").Value = '''56.06'
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("D39").Value = '''2.68'
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("E40").Value = '  -3.66%  '
$ws.Range("D41").Value = '0.0₃0697'
$ws.Range("E41").Value = '  -5.13%  '
$ws.Range("D42").Value = '''32.79'
$ws.Range("E42").Value = '  -5.03%  '
$ws.Range("D43").Value = '''3.44'
$ws.Range("E43").Value = '  +5.91%  '
$ws.Range("D44").Value = '''3.18'
$ws.Range("E44").Value = '  -4.60%  '
$ws.Range("D45").Value = '''0.336'
$ws.Range("E45").Value = '  -2.60%  '
$ws.Range("E46").Value = '  -3.38%  '
$ws.Range("D47").Value = '''3.02'
$ws.Range("E47").Value = '  -11.74%  '
$ws.Range("E48").Value = '  -2.15%  '
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("E50").Value = '  -3.20%  '
$ws.Range("E51").Value = '  +1.22%  '
